$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(2, 1).Range.Text = "guest_gregs"
$t.Cell(3, 1).Range.Text = "guest_pfreitas"
$t.Cell(4, 1).Range.Text = "guest_sono"
$t.Cell(5, 1).Range.Text = "nsilva"
$t.Cell(6, 1).Range.Text = "rpinheiro"
$t.Cell(7, 1).Range.Text = "ggomes"
$t.Cell(8, 1).Range.Text = "guest_diogo"
$t.Cell(9, 1).Range.Text = "guest_cgomes"
$t.Cell(11, 1).Range.Text = "pduarte"
